$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 127
$ws1.Range("F5").Value = 1746
$ws1.Range("F6").Value = 3316
$ws1.Range("F7").Value = 1003
$ws1.Range("F8").Value = 2172
$ws1.Range("F9").Value = 2089
$ws1.Range("F10").Value = 1088
$ws1.Range("F11").Value = 592
$ws1.Range("F12").Value = 17
$ws1.Range("F13").Value = 1655
$ws1.Range("F14").Value = 383
$ws1.Range("F16").Value = 36
$ws1.Range("F17").Value = 89
$ws1.Range("F18").Value = 187
$ws1.Range("F19").Value = 1553
$ws1.Range("F20").Value = 602
$ws1.Range("F21").Value = 703
$ws1.Range("F22").Value = 583
$ws1.Range("F23").Value = 12144
$ws1.Range("F24").Value = 12169
$ws1.Range("F25").Value = 902
$ws1.Range("F26").Value = 688
$ws1.Range("F28").Value = 21
$ws1.Range("F29").Value = 14
$ws1.Range("F30").Value = 330
$ws1.Range("F31").Value = 1908
$ws1.Range("F32").Value = 188
$ws1.Range("F33").Value = 560

# Sheet "全部类型" (sheet4.xml) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 127
$ws4.Range("F6").Value = 1746
$ws4.Range("F7").Value = 3316
$ws4.Range("F8").Value = 1005
$ws4.Range("F9").Value = 2172
$ws4.Range("F10").Value = 2089
$ws4.Range("F11").Value = 1088
$ws4.Range("F12").Value = 592
$ws4.Range("F13").Value = 17
$ws4.Range("F14").Value = 1655
$ws4.Range("F15").Value = 383
$ws4.Range("F18").Value = 36
$ws4.Range("F20").Value = 89
$ws4.Range("F22").Value = 187
$ws4.Range("F23").Value = 1553
$ws4.Range("F24").Value = 602
$ws4.Range("F25").Value = 703
$ws4.Range("F26").Value = 583
$ws4.Range("F27").Value = 12144
$ws4.Range("F28").Value = 12169
$ws4.Range("F29").Value = 902
$ws4.Range("F30").Value = 688
$ws4.Range("F32").Value = 21
$ws4.Range("F33").Value = 14
$ws4.Range("F34").Value = 330
$ws4.Range("F35").Value = 1908
$ws4.Range("F38").Value = 188
$ws4.Range("F39").Value = 560
